$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$arr = New-Object 'object[,]' 96,2

$arr[0,0] = 45740.01041666666
$arr[0,1] = 0
$arr[1,0] = 45740.02083333334
$arr[1,1] = 0
$arr[2,0] = 45740.03125
$arr[2,1] = 0
$arr[3,0] = 45740.04166666666
$arr[3,1] = 0
$arr[4,0] = 45740.05208333334
$arr[4,1] = 0
$arr[5,0] = 45740.0625
$arr[5,1] = 0
$arr[6,0] = 45740.07291666666
$arr[6,1] = 0
$arr[7,0] = 45740.08333333334
$arr[7,1] = 0
$arr[8,0] = 45740.09375
$arr[8,1] = 0
$arr[9,0] = 45740.10416666666
$arr[9,1] = 0
$arr[10,0] = 45740.11458333334
$arr[10,1] = 0
$arr[11,0] = 45740.125
$arr[11,1] = 0
$arr[12,0] = 45740.13541666666
$arr[12,1] = 0
$arr[13,0] = 45740.14583333334
$arr[13,1] = 0
$arr[14,0] = 45740.15625
$arr[14,1] = 0
$arr[15,0] = 45740.16666666666
$arr[15,1] = 0
$arr[16,0] = 45740.17708333334
$arr[16,1] = 0
$arr[17,0] = 45740.1875
$arr[17,1] = 0
$arr[18,0] = 45740.19791666666
$arr[18,1] = 0
$arr[19,0] = 45740.20833333334
$arr[19,1] = 0
$arr[20,0] = 45740.21875
$arr[20,1] = 20
$arr[21,0] = 45740.22916666666
$arr[21,1] = 17
$arr[22,0] = 45740.23958333334
$arr[22,1] = 18
$arr[23,0] = 45740.25
$arr[23,1] = 22
$arr[24,0] = 45740.26041666666
$arr[24,1] = 182
$arr[25,0] = 45740.27083333334
$arr[25,1] = 195
$arr[26,0] = 45740.28125
$arr[26,1] = 211
$arr[27,0] = 45740.29166666666
$arr[27,1] = 234
$arr[28,0] = 45740.30208333334
$arr[28,1] = 487
$arr[29,0] = 45740.3125
$arr[29,1] = 540
$arr[30,0] = 45740.32291666666
$arr[30,1] = 567
$arr[31,0] = 45740.33333333334
$arr[31,1] = 597
$arr[32,0] = 45740.34375
$arr[32,1] = 825
$arr[33,0] = 45740.35416666666
$arr[33,1] = 855
$arr[34,0] = 45740.36458333334
$arr[34,1] = 883
$arr[35,0] = 45740.375
$arr[35,1] = 911
$arr[36,0] = 45740.38541666666
$arr[36,1] = 1098
$arr[37,0] = 45740.39583333334
$arr[37,1] = 1118
$arr[38,0] = 45740.40625
$arr[38,1] = 1136
$arr[39,0] = 45740.41666666666
$arr[39,1] = 1151
$arr[40,0] = 45740.42708333334
$arr[40,1] = 1215
$arr[41,0] = 45740.4375
$arr[41,1] = 1224
$arr[42,0] = 45740.44791666666
$arr[42,1] = 1228
$arr[43,0] = 45740.45833333334
$arr[43,1] = 1232
$arr[44,0] = 45740.46875
$arr[44,1] = 1204
$arr[45,0] = 45740.47916666666
$arr[45,1] = 1201
$arr[46,0] = 45740.48958333334
$arr[46,1] = 1197
$arr[47,0] = 45740.5
$arr[47,1] = 1190
$arr[48,0] = 45740.51041666666
$arr[48,1] = 1105
$arr[49,0] = 45740.52083333334
$arr[49,1] = 1096
$arr[50,0] = 45740.53125
$arr[50,1] = 1083
$arr[51,0] = 45740.54166666666
$arr[51,1] = 1068
$arr[52,0] = 45740.55208333334
$arr[52,1] = 935
$arr[53,0] = 45740.5625
$arr[53,1] = 916
$arr[54,0] = 45740.57291666666
$arr[54,1] = 897
$arr[55,0] = 45740.58333333334
$arr[55,1] = 874
$arr[56,0] = 45740.59375
$arr[56,1] = 713
$arr[57,0] = 45740.60416666666
$arr[57,1] = 692
$arr[58,0] = 45740.61458333334
$arr[58,1] = 667
$arr[59,0] = 45740.625
$arr[59,1] = 644
$arr[60,0] = 45740.63541666666
$arr[60,1] = 431
$arr[61,0] = 45740.64583333334
$arr[61,1] = 409
$arr[62,0] = 45740.65625
$arr[62,1] = 390
$arr[63,0] = 45740.66666666666
$arr[63,1] = 372
$arr[64,0] = 45740.67708333334
$arr[64,1] = 181
$arr[65,0] = 45740.6875
$arr[65,1] = 168
$arr[66,0] = 45740.69791666666
$arr[66,1] = 153
$arr[67,0] = 45740.70833333334
$arr[67,1] = 141
$arr[68,0] = 45740.71875
$arr[68,1] = 16
$arr[69,0] = 45740.72916666666
$arr[69,1] = 12
$arr[70,0] = 45740.73958333334
$arr[70,1] = 10
$arr[71,0] = 45740.75
$arr[71,1] = 8
$arr[72,0] = 45740.76041666666
$arr[72,1] = 1
$arr[73,0] = 45740.77083333334
$arr[73,1] = 1
$arr[74,0] = 45740.78125
$arr[74,1] = 0
$arr[75,0] = 45740.79166666666
$arr[75,1] = 0
$arr[76,0] = 45740.80208333334
$arr[76,1] = 0
$arr[77,0] = 45740.8125
$arr[77,1] = 0
$arr[78,0] = 45740.82291666666
$arr[78,1] = 0
$arr[79,0] = 45740.83333333334
$arr[79,1] = 0
$arr[80,0] = 45740.84375
$arr[80,1] = 0
$arr[81,0] = 45740.85416666666
$arr[81,1] = 0
$arr[82,0] = 45740.86458333334
$arr[82,1] = 0
$arr[83,0] = 45740.875
$arr[83,1] = 0
$arr[84,0] = 45740.88541666666
$arr[84,1] = 0
$arr[85,0] = 45740.89583333334
$arr[85,1] = 0
$arr[86,0] = 45740.90625
$arr[86,1] = 0
$arr[87,0] = 45740.91666666666
$arr[87,1] = 0
$arr[88,0] = 45740.92708333334
$arr[88,1] = 0
$arr[89,0] = 45740.9375
$arr[89,1] = 0
$arr[90,0] = 45740.94791666666
$arr[90,1] = 0
$arr[91,0] = 45740.95833333334
$arr[91,1] = 0
$arr[92,0] = 45740.96875
$arr[92,1] = 0
$arr[93,0] = 45740.97916666666
$arr[93,1] = 0
$arr[94,0] = 45740.98958333334
$arr[94,1] = 0
$arr[95,0] = 45741
$arr[95,1] = 0

$ws.Range("A2:B97").Value = $arr
